# Apply updates to the "展览" and "全部类型" sheets:
#  - F2: 118 -> 120
#  - C3: "丽水·LPJ 现实X次元动漫展" -> "丽水·LPJ 现实X次元动漫展（取消）"
#  - G3: 45 (number) -> "不可售" (text, event is no longer sellable)

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 120
    $ws.Range("C3").Value = "丽水·LPJ 现实X次元动漫展（取消）"
    $ws.Range("G3").Value = "不可售"
}
